$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("titers")

# Row 167 is an exact duplicate of row 166 (same label, same measurements for
# most columns, with only O167 differing) -- remove the duplicate row and
# shift the remaining data rows up, as described in the commit message
# ("Removed duplicate row from sample data").
$ws.Rows.Item(167).Delete()

# Reflect the new selection: the (now shifted) row 167 is selected in full.
$ws.Activate()
$ws.Range("A167:XFD167").Select()
